$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "From DataTable" sample: normalize the Date column (D9:D13) ---
# Previously every row held the same timestamp (6/26/2011 6:32 AM). The
# refreshed example now uses plain, incrementing calendar dates with no
# time-of-day component, and a date-only display format.
$dates = @("2000-01-01", "2000-01-02", "2000-01-03", "2000-01-04", "2000-01-05")
$row = 9
foreach ($d in $dates) {
    $ws.Cells.Item($row, 4).Value = [DateTime]$d
    $row++
}

# Date-only format (standard/built-in format id 14, displayed as m/d/yyyy)
# instead of the previous date+time format (built-in id 22).
$ws.Range("D9:D13").NumberFormat = "mm-dd-yy"

# --- Column widths: column D now matches the width of its C/E neighbours ---
$ws.Columns("D:D").ColumnWidth = $ws.Columns("C:C").ColumnWidth
